$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 0) Fix CodigoNiv2 (column J) text codes for rows 169-173 to proper 4-digit codes
$ws.Cells.Item(169, 10).Value = '0406'
$ws.Cells.Item(170, 10).Value = '0408'
$ws.Cells.Item(171, 10).Value = '0419'
$ws.Cells.Item(172, 10).Value = '0413'
$ws.Cells.Item(173, 10).Value = '0410'

# 2) Extend the HOSPITALES_HN table with 21 new rows (174-194)
$lo = $ws.ListObjects.Item(1)
for ($i = 0; $i -lt 21; $i++) {
    $null = $lo.ListRows.Add()
}

# 3) Populate the new rows with data
# row 174
$ws.Cells.Item(174, 2).Value = 'HND'
$ws.Cells.Item(174, 3).Value = 'Honduras'
$ws.Cells.Item(174, 4).Value = 3
$ws.Cells.Item(174, 5).Value = 6
$ws.Cells.Item(174, 6).Value = '06'
$ws.Cells.Item(174, 7).Value = 'Choluteca'
$ws.Cells.Item(174, 8).Value = 'Departamento'
$ws.Cells.Item(174, 9).Value = 12
$ws.Cells.Item(174, 10).Value = '0612'
$ws.Cells.Item(174, 11).Value = 'San Antonio de Flores'
$ws.Cells.Item(174, 12).Value = 'Municipio'
$ws.Cells.Item(174, 13).Value = 1
$ws.Cells.Item(174, 14).Value = '061201'
$ws.Cells.Item(174, 15).Value = 'San Antonio de Flores'
$ws.Cells.Item(174, 16).Value = 'Aldea'
$ws.Cells.Item(174, 17).Value = 'HND-0612'
$ws.Cells.Item(174, 18).Value = 'Salud'
$ws.Cells.Item(174, 19).Value = 'Cesamo'
$ws.Cells.Item(174, 22).Value = 13.665829
$ws.Cells.Item(174, 23).Value = -87.362112
# row 175
$ws.Cells.Item(175, 2).Value = 'HND'
$ws.Cells.Item(175, 3).Value = 'Honduras'
$ws.Cells.Item(175, 4).Value = 3
$ws.Cells.Item(175, 5).Value = 6
$ws.Cells.Item(175, 6).Value = '06'
$ws.Cells.Item(175, 7).Value = 'Choluteca'
$ws.Cells.Item(175, 8).Value = 'Departamento'
$ws.Cells.Item(175, 9).Value = 8
$ws.Cells.Item(175, 10).Value = '0608'
$ws.Cells.Item(175, 11).Value = 'Morolica'
$ws.Cells.Item(175, 12).Value = 'Municipio'
$ws.Cells.Item(175, 13).Value = 1
$ws.Cells.Item(175, 14).Value = '060801'
$ws.Cells.Item(175, 15).Value = 'Morolica'
$ws.Cells.Item(175, 16).Value = 'Aldea'
$ws.Cells.Item(175, 17).Value = 'HND-0608'
$ws.Cells.Item(175, 18).Value = 'Salud'
$ws.Cells.Item(175, 19).Value = 'Cesamo'
$ws.Cells.Item(175, 22).Value = 13.568473
$ws.Cells.Item(175, 23).Value = -86.907604
# row 176
$ws.Cells.Item(176, 2).Value = 'HND'
$ws.Cells.Item(176, 3).Value = 'Honduras'
$ws.Cells.Item(176, 4).Value = 3
$ws.Cells.Item(176, 5).Value = 6
$ws.Cells.Item(176, 6).Value = '06'
$ws.Cells.Item(176, 7).Value = 'Choluteca'
$ws.Cells.Item(176, 8).Value = 'Departamento'
$ws.Cells.Item(176, 9).Value = 16
$ws.Cells.Item(176, 10).Value = '0616'
$ws.Cells.Item(176, 11).Value = 'Santa Ana de Yusguare'
$ws.Cells.Item(176, 12).Value = 'Municipio'
$ws.Cells.Item(176, 13).Value = 1
$ws.Cells.Item(176, 14).Value = '061601'
$ws.Cells.Item(176, 15).Value = 'Santa Ana de Yusguare'
$ws.Cells.Item(176, 16).Value = 'Aldea'
$ws.Cells.Item(176, 17).Value = 'HND-0616'
$ws.Cells.Item(176, 18).Value = 'Salud'
$ws.Cells.Item(176, 19).Value = 'Cesamo'
$ws.Cells.Item(176, 22).Value = 13.293342
$ws.Cells.Item(176, 23).Value = -87.111866
# row 177
$ws.Cells.Item(177, 2).Value = 'HND'
$ws.Cells.Item(177, 3).Value = 'Honduras'
$ws.Cells.Item(177, 4).Value = 3
$ws.Cells.Item(177, 5).Value = 6
$ws.Cells.Item(177, 6).Value = '06'
$ws.Cells.Item(177, 7).Value = 'Choluteca'
$ws.Cells.Item(177, 8).Value = 'Departamento'
$ws.Cells.Item(177, 9).Value = 7
$ws.Cells.Item(177, 10).Value = '0607'
$ws.Cells.Item(177, 11).Value = 'Marcovia'
$ws.Cells.Item(177, 12).Value = 'Municipio'
$ws.Cells.Item(177, 13).Value = 1
$ws.Cells.Item(177, 14).Value = '060701'
$ws.Cells.Item(177, 15).Value = 'Marcovia'
$ws.Cells.Item(177, 16).Value = 'Aldea'
$ws.Cells.Item(177, 17).Value = 'HND-0607'
$ws.Cells.Item(177, 18).Value = 'Salud'
$ws.Cells.Item(177, 19).Value = 'Cesamo'
$ws.Cells.Item(177, 22).Value = 13.284372
$ws.Cells.Item(177, 23).Value = -87.312592
# row 178
$ws.Cells.Item(178, 2).Value = 'HND'
$ws.Cells.Item(178, 3).Value = 'Honduras'
$ws.Cells.Item(178, 4).Value = 3
$ws.Cells.Item(178, 5).Value = 6
$ws.Cells.Item(178, 6).Value = '06'
$ws.Cells.Item(178, 7).Value = 'Choluteca'
$ws.Cells.Item(178, 8).Value = 'Departamento'
$ws.Cells.Item(178, 9).Value = 9
$ws.Cells.Item(178, 10).Value = '0609'
$ws.Cells.Item(178, 11).Value = 'Namasigue'
$ws.Cells.Item(178, 12).Value = 'Municipio'
$ws.Cells.Item(178, 13).Value = 1
$ws.Cells.Item(178, 14).Value = '060901'
$ws.Cells.Item(178, 15).Value = 'Namasigue'
$ws.Cells.Item(178, 16).Value = 'Aldea'
$ws.Cells.Item(178, 17).Value = 'HND-0609'
$ws.Cells.Item(178, 18).Value = 'Salud'
$ws.Cells.Item(178, 19).Value = 'Cesamo'
$ws.Cells.Item(178, 22).Value = 13.203086
$ws.Cells.Item(178, 23).Value = -87.13876
# row 179
$ws.Cells.Item(179, 2).Value = 'HND'
$ws.Cells.Item(179, 3).Value = 'Honduras'
$ws.Cells.Item(179, 4).Value = 3
$ws.Cells.Item(179, 5).Value = 6
$ws.Cells.Item(179, 6).Value = '06'
$ws.Cells.Item(179, 7).Value = 'Choluteca'
$ws.Cells.Item(179, 8).Value = 'Departamento'
$ws.Cells.Item(179, 9).Value = 7
$ws.Cells.Item(179, 10).Value = '0607'
$ws.Cells.Item(179, 11).Value = 'Marcovia'
$ws.Cells.Item(179, 12).Value = 'Municipio'
$ws.Cells.Item(179, 13).Value = 14
$ws.Cells.Item(179, 14).Value = '060714'
$ws.Cells.Item(179, 15).Value = 'Monjaras'
$ws.Cells.Item(179, 16).Value = 'Aldea'
$ws.Cells.Item(179, 17).Value = 'HND-0607'
$ws.Cells.Item(179, 18).Value = 'Salud'
$ws.Cells.Item(179, 19).Value = 'Cesamo'
$ws.Cells.Item(179, 22).Value = 13.198708
$ws.Cells.Item(179, 23).Value = -87.375153
# row 180
$ws.Cells.Item(180, 2).Value = 'HND'
$ws.Cells.Item(180, 3).Value = 'Honduras'
$ws.Cells.Item(180, 4).Value = 3
$ws.Cells.Item(180, 5).Value = 8
$ws.Cells.Item(180, 6).Value = '08'
$ws.Cells.Item(180, 7).Value = 'Francisco Morazán'
$ws.Cells.Item(180, 8).Value = 'Departamento'
$ws.Cells.Item(180, 9).Value = 13
$ws.Cells.Item(180, 10).Value = '0813'
$ws.Cells.Item(180, 11).Value = 'Ojojona'
$ws.Cells.Item(180, 12).Value = 'Municipio'
$ws.Cells.Item(180, 15).Value = 'Ojojona'
$ws.Cells.Item(180, 16).Value = 'Aldea'
$ws.Cells.Item(180, 17).Value = 'HND-0813'
$ws.Cells.Item(180, 18).Value = 'Salud'
$ws.Cells.Item(180, 19).Value = 'Cesamo'
$ws.Cells.Item(180, 22).Value = 13.932455
$ws.Cells.Item(180, 23).Value = -87.297745
# row 181
$ws.Cells.Item(181, 2).Value = 'HND'
$ws.Cells.Item(181, 3).Value = 'Honduras'
$ws.Cells.Item(181, 4).Value = 3
$ws.Cells.Item(181, 5).Value = 8
$ws.Cells.Item(181, 6).Value = '08'
$ws.Cells.Item(181, 7).Value = 'Francisco Morazán'
$ws.Cells.Item(181, 8).Value = 'Departamento'
$ws.Cells.Item(181, 9).Value = 11
$ws.Cells.Item(181, 10).Value = '0811'
$ws.Cells.Item(181, 11).Value = 'Marale'
$ws.Cells.Item(181, 12).Value = 'Municipio'
$ws.Cells.Item(181, 13).Value = 2
$ws.Cells.Item(181, 14).Value = '081102'
$ws.Cells.Item(181, 15).Value = 'Marale'
$ws.Cells.Item(181, 16).Value = 'Aldea'
$ws.Cells.Item(181, 17).Value = 'HND-0811'
$ws.Cells.Item(181, 18).Value = 'Salud'
$ws.Cells.Item(181, 19).Value = 'Cesamo'
$ws.Cells.Item(181, 22).Value = 14.915878
$ws.Cells.Item(181, 23).Value = -87.194402
# row 182
$ws.Cells.Item(182, 2).Value = 'HND'
$ws.Cells.Item(182, 3).Value = 'Honduras'
$ws.Cells.Item(182, 4).Value = 3
$ws.Cells.Item(182, 5).Value = 8
$ws.Cells.Item(182, 6).Value = '08'
$ws.Cells.Item(182, 7).Value = 'Francisco Morazán'
$ws.Cells.Item(182, 8).Value = 'Departamento'
$ws.Cells.Item(182, 9).Value = 5
$ws.Cells.Item(182, 10).Value = '0805'
$ws.Cells.Item(182, 11).Value = 'El Porvenir'
$ws.Cells.Item(182, 12).Value = 'Municipio'
$ws.Cells.Item(182, 13).Value = 1
$ws.Cells.Item(182, 14).Value = '080501'
$ws.Cells.Item(182, 15).Value = 'El Porvenir'
$ws.Cells.Item(182, 16).Value = 'Aldea'
$ws.Cells.Item(182, 17).Value = 'HND-0805'
$ws.Cells.Item(182, 18).Value = 'Salud'
$ws.Cells.Item(182, 19).Value = 'Cesamo'
$ws.Cells.Item(182, 22).Value = 14.754858
$ws.Cells.Item(182, 23).Value = -87.188081
# row 183
$ws.Cells.Item(183, 2).Value = 'HND'
$ws.Cells.Item(183, 3).Value = 'Honduras'
$ws.Cells.Item(183, 4).Value = 3
$ws.Cells.Item(183, 5).Value = 8
$ws.Cells.Item(183, 6).Value = '08'
$ws.Cells.Item(183, 7).Value = 'Francisco Morazán'
$ws.Cells.Item(183, 8).Value = 'Departamento'
$ws.Cells.Item(183, 9).Value = 14
$ws.Cells.Item(183, 10).Value = '0814'
$ws.Cells.Item(183, 11).Value = 'Orica'
$ws.Cells.Item(183, 12).Value = 'Municipio'
$ws.Cells.Item(183, 13).Value = 1
$ws.Cells.Item(183, 14).Value = '081401'
$ws.Cells.Item(183, 15).Value = 'Orica'
$ws.Cells.Item(183, 16).Value = 'Aldea'
$ws.Cells.Item(183, 17).Value = 'HND-0814'
$ws.Cells.Item(183, 18).Value = 'Salud'
$ws.Cells.Item(183, 19).Value = 'Cesamo'
$ws.Cells.Item(183, 22).Value = 14.714986
$ws.Cells.Item(183, 23).Value = -86.942894
# row 184
$ws.Cells.Item(184, 2).Value = 'HND'
$ws.Cells.Item(184, 3).Value = 'Honduras'
$ws.Cells.Item(184, 4).Value = 3
$ws.Cells.Item(184, 5).Value = 8
$ws.Cells.Item(184, 6).Value = '08'
$ws.Cells.Item(184, 7).Value = 'Francisco Morazán'
$ws.Cells.Item(184, 8).Value = 'Departamento'
$ws.Cells.Item(184, 9).Value = 19
$ws.Cells.Item(184, 10).Value = '0819'
$ws.Cells.Item(184, 11).Value = 'San Ignacio'
$ws.Cells.Item(184, 12).Value = 'Municipio'
$ws.Cells.Item(184, 13).Value = 1
$ws.Cells.Item(184, 15).Value = 'San Ignacio'
$ws.Cells.Item(184, 16).Value = 'Aldea'
$ws.Cells.Item(184, 17).Value = 'HND-'
$ws.Cells.Item(184, 18).Value = 'Salud'
$ws.Cells.Item(184, 19).Value = 'Cesamo'
# row 185
$ws.Cells.Item(185, 2).Value = 'HND'
$ws.Cells.Item(185, 3).Value = 'Honduras'
$ws.Cells.Item(185, 4).Value = 3
$ws.Cells.Item(185, 8).Value = 'Departamento'
$ws.Cells.Item(185, 12).Value = 'Municipio'
$ws.Cells.Item(185, 16).Value = 'Aldea'
$ws.Cells.Item(185, 17).Value = 'HND-'
$ws.Cells.Item(185, 18).Value = 'Salud'
$ws.Cells.Item(185, 19).Value = 'Cesamo'
# row 186
$ws.Cells.Item(186, 2).Value = 'HND'
$ws.Cells.Item(186, 3).Value = 'Honduras'
$ws.Cells.Item(186, 4).Value = 3
$ws.Cells.Item(186, 8).Value = 'Departamento'
$ws.Cells.Item(186, 12).Value = 'Municipio'
$ws.Cells.Item(186, 16).Value = 'Aldea'
$ws.Cells.Item(186, 17).Value = 'HND-'
$ws.Cells.Item(186, 18).Value = 'Salud'
$ws.Cells.Item(186, 19).Value = 'Cesamo'
# row 187
$ws.Cells.Item(187, 2).Value = 'HND'
$ws.Cells.Item(187, 3).Value = 'Honduras'
$ws.Cells.Item(187, 4).Value = 3
$ws.Cells.Item(187, 8).Value = 'Departamento'
$ws.Cells.Item(187, 12).Value = 'Municipio'
$ws.Cells.Item(187, 16).Value = 'Aldea'
$ws.Cells.Item(187, 17).Value = 'HND-'
$ws.Cells.Item(187, 18).Value = 'Salud'
$ws.Cells.Item(187, 19).Value = 'Cesamo'
# row 188
$ws.Cells.Item(188, 2).Value = 'HND'
$ws.Cells.Item(188, 3).Value = 'Honduras'
$ws.Cells.Item(188, 4).Value = 3
$ws.Cells.Item(188, 8).Value = 'Departamento'
$ws.Cells.Item(188, 12).Value = 'Municipio'
$ws.Cells.Item(188, 16).Value = 'Aldea'
$ws.Cells.Item(188, 17).Value = 'HND-'
$ws.Cells.Item(188, 18).Value = 'Salud'
$ws.Cells.Item(188, 19).Value = 'Cesamo'
# row 189
$ws.Cells.Item(189, 2).Value = 'HND'
$ws.Cells.Item(189, 3).Value = 'Honduras'
$ws.Cells.Item(189, 4).Value = 3
$ws.Cells.Item(189, 8).Value = 'Departamento'
$ws.Cells.Item(189, 12).Value = 'Municipio'
$ws.Cells.Item(189, 16).Value = 'Aldea'
$ws.Cells.Item(189, 17).Value = 'HND-'
$ws.Cells.Item(189, 18).Value = 'Salud'
$ws.Cells.Item(189, 19).Value = 'Cesamo'
# row 190
$ws.Cells.Item(190, 2).Value = 'HND'
$ws.Cells.Item(190, 3).Value = 'Honduras'
$ws.Cells.Item(190, 4).Value = 3
$ws.Cells.Item(190, 8).Value = 'Departamento'
$ws.Cells.Item(190, 12).Value = 'Municipio'
$ws.Cells.Item(190, 16).Value = 'Aldea'
$ws.Cells.Item(190, 17).Value = 'HND-'
$ws.Cells.Item(190, 18).Value = 'Salud'
$ws.Cells.Item(190, 19).Value = 'Cesamo'
# row 191
$ws.Cells.Item(191, 2).Value = 'HND'
$ws.Cells.Item(191, 3).Value = 'Honduras'
$ws.Cells.Item(191, 4).Value = 3
$ws.Cells.Item(191, 8).Value = 'Departamento'
$ws.Cells.Item(191, 12).Value = 'Municipio'
$ws.Cells.Item(191, 16).Value = 'Aldea'
$ws.Cells.Item(191, 17).Value = 'HND-'
$ws.Cells.Item(191, 18).Value = 'Salud'
$ws.Cells.Item(191, 19).Value = 'Cesamo'
# row 192
$ws.Cells.Item(192, 2).Value = 'HND'
$ws.Cells.Item(192, 3).Value = 'Honduras'
$ws.Cells.Item(192, 4).Value = 3
$ws.Cells.Item(192, 8).Value = 'Departamento'
$ws.Cells.Item(192, 12).Value = 'Municipio'
$ws.Cells.Item(192, 16).Value = 'Aldea'
$ws.Cells.Item(192, 17).Value = 'HND-'
$ws.Cells.Item(192, 18).Value = 'Salud'
$ws.Cells.Item(192, 19).Value = 'Cesamo'
# row 193
$ws.Cells.Item(193, 2).Value = 'HND'
$ws.Cells.Item(193, 3).Value = 'Honduras'
$ws.Cells.Item(193, 4).Value = 3
$ws.Cells.Item(193, 8).Value = 'Departamento'
$ws.Cells.Item(193, 12).Value = 'Municipio'
$ws.Cells.Item(193, 16).Value = 'Aldea'
$ws.Cells.Item(193, 17).Value = 'HND-'
$ws.Cells.Item(193, 18).Value = 'Salud'
$ws.Cells.Item(193, 19).Value = 'Cesamo'
# row 194
$ws.Cells.Item(194, 2).Value = 'HND'
$ws.Cells.Item(194, 3).Value = 'Honduras'
$ws.Cells.Item(194, 4).Value = 3
$ws.Cells.Item(194, 8).Value = 'Departamento'
$ws.Cells.Item(194, 12).Value = 'Municipio'
$ws.Cells.Item(194, 16).Value = 'Aldea'
$ws.Cells.Item(194, 17).Value = 'HND-'
$ws.Cells.Item(194, 18).Value = 'Salud'
$ws.Cells.Item(194, 19).Value = 'Cesamo'

# 4) Update the _FilterDatabase defined name to cover the new table range
foreach ($n in $wb.Names) {
    if ($n.Name -eq "HOSPITALES!_FilterDatabase") {
        $n.RefersTo = "=HOSPITALES!`$A`$1:`$W`$194"
    }
}

